$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 698, pushing the existing rows 698-737 down to 701-740.
$ws.Range("A698:T700").EntireRow.Insert()

# Populate the 3 newly inserted rows with the new weekly price data
# (same market / product metadata as every other row in this sheet).

# Row 698 - Especial
$ws.Cells.Item(698, 1).Value = 8
$ws.Cells.Item(698, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(698, 3).Value = "Coquimbo"
$ws.Cells.Item(698, 4).Value = 44706
$ws.Cells.Item(698, 5).Value = 4
$ws.Cells.Item(698, 6).Value = "Fruta"
$ws.Cells.Item(698, 7).Value = 100101
$ws.Cells.Item(698, 8).Value = "Berries"
$ws.Cells.Item(698, 9).Value = 100112025
$ws.Cells.Item(698, 10).Value = "Frutilla"
$ws.Cells.Item(698, 11).Value = "Sin especificar"
$ws.Cells.Item(698, 12).Value = "Especial"
$ws.Cells.Item(698, 13).Value = 300
$ws.Cells.Item(698, 14).Value = 18000
$ws.Cells.Item(698, 15).Value = 19000
$ws.Cells.Item(698, 16).Value = 18500
$ws.Cells.Item(698, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(698, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(698, 19).Value = 2643
$ws.Cells.Item(698, 20).Value = 7

# Row 699 - Primera
$ws.Cells.Item(699, 1).Value = 8
$ws.Cells.Item(699, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(699, 3).Value = "Coquimbo"
$ws.Cells.Item(699, 4).Value = 44706
$ws.Cells.Item(699, 5).Value = 4
$ws.Cells.Item(699, 6).Value = "Fruta"
$ws.Cells.Item(699, 7).Value = 100101
$ws.Cells.Item(699, 8).Value = "Berries"
$ws.Cells.Item(699, 9).Value = 100112025
$ws.Cells.Item(699, 10).Value = "Frutilla"
$ws.Cells.Item(699, 11).Value = "Sin especificar"
$ws.Cells.Item(699, 12).Value = "Primera"
$ws.Cells.Item(699, 13).Value = 400
$ws.Cells.Item(699, 14).Value = 16000
$ws.Cells.Item(699, 15).Value = 17000
$ws.Cells.Item(699, 16).Value = 16500
$ws.Cells.Item(699, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(699, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(699, 19).Value = 2357
$ws.Cells.Item(699, 20).Value = 7

# Row 700 - Segunda
$ws.Cells.Item(700, 1).Value = 8
$ws.Cells.Item(700, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(700, 3).Value = "Coquimbo"
$ws.Cells.Item(700, 4).Value = 44706
$ws.Cells.Item(700, 5).Value = 4
$ws.Cells.Item(700, 6).Value = "Fruta"
$ws.Cells.Item(700, 7).Value = 100101
$ws.Cells.Item(700, 8).Value = "Berries"
$ws.Cells.Item(700, 9).Value = 100112025
$ws.Cells.Item(700, 10).Value = "Frutilla"
$ws.Cells.Item(700, 11).Value = "Sin especificar"
$ws.Cells.Item(700, 12).Value = "Segunda"
$ws.Cells.Item(700, 13).Value = 270
$ws.Cells.Item(700, 14).Value = 12000
$ws.Cells.Item(700, 15).Value = 13000
$ws.Cells.Item(700, 16).Value = 12444
$ws.Cells.Item(700, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(700, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(700, 19).Value = 1778
$ws.Cells.Item(700, 20).Value = 7
